$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13, pushing "Programa resumido:" and everything below
# down by one (row 13 becomes blank, to be used for the Docentes responsaveis
# name that used to incorrectly live under "Objetivos:").
$ws.Rows("13:13").Insert()

# --- Objetivos: (row 10) -------------------------------------------------
# Replace the mis-placed professor name with the real objectives text.
$objetivos = "Formação dos estudantes de Engenharia Bioquímica na área de tecnologia de enzimas, com foco principal nos estudos de estrutura versus propriedades e mecanismos de ação, controle operacional na purificação e imobilização de enzimas, formas de determinação de atividade enzimática e aplicações das enzimas nos processos industriais."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# --- Docentes responsaveis: (row 12) now has its data on the new row 13 --
$ws.Range("B13").Value = "4873328 - Fernando Segato"
$ws.Range("C13").Value = "4873328 - Fernando Segato"

# --- Programa resumido: (row 14) -----------------------------------------
$programaResumido = "A disciplina aborda como as enzimas atuam, como se definem as estratégias de purificação e quais são as principais aplicações tecnológicas das enzimas. Dentro dos processos de purificação, o foco envolve a definição de estratégias apropriadas para a purificação em etapas sequenciais, os métodos de controle de cada etapa, além dos métodos de monitoramento da atividade enzimática. Também se aborda a aplicação das enzimas em processos industriais."
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Row 15 "Short syllabus:" already carries the correct English text after
# the row shift, no change needed there.

# --- Programa: (row 16) ---------------------------------------------------
$programa = "1. Origem celular das enzimas: origem das enzimas, diferenciação entre enzimas intra e extracelulares, importância fisiológica e introdução ao mercado mundial de enzimas.2. Estrutura versus propriedades e mecanismos de ação das enzimas: estruturas tridimensionais e sua determinação, importância da estrutura terciária na atividade catalítica, ação catalítica de proteases, glicosidases e oxido-redutases.3. Controle operacional na purificação de enzimas: métodos de extração de enzimas, métodos de purificação preliminar, métodos de separação baseados na carga, no tamanho e na afinidade. Definição de estratégias de purificação de enzimas.4. Métodos de determinação de atividade enzimática: definição de atividade em unidades internacionais, significado da atividade enzimática, formas de quantificar e expressar e atividade enzimática. Requerimentos de um método experimental usado na determinação de atividade enzimática.5. Cinética enzimática: métodos gráficos e numéricos de determinação de velocidade inicial de reação, condições experimentais demandadas para determinar a velocidade inicial, cálculos de atividade enzimática.6. Enzimas imobilizadas: formas de imobilização e aplicações de sistemas imobilizados.7. Aplicações de enzimas na indústria: uso de enzimas em detergentes, no processamento do amido, na indústria alimentícia, na indústria têxtil, na síntese de fármacos e na indústria de celulose e papel."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Row 17 "Syllabus:" already carries the correct English text after the
# row shift, no change needed there.

# --- Método: (row 19) ------------------------------------------------------
$metodo = "A avaliação será feita por meio de provas escritas (P1 e P2)."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Critério: (row 20) -----------------------------------------------------
$criterio = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1x1 + P2x2)/3"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# --- Norma de recuperação: (row 21) -----------------------------------------
$norma = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# --- Bibliografia: (row 22) --------------------------------------------------
$bibliografia = "1. BON, E.S., FERRARA M.A., CORVO M.L. (Eds.) Enzimas em Biotecnologia - Produção, aplicação e mercado, Rio de Janeiro: Editora Interciêcnia, 2008.`n2. COPELAND, R.A. Enzymes: a practical introduction to structure, mechanism and data analysis, New York: Academic Press, 2000.`n3. LEHNINGER, A.L., NELSON, O.L., COX, M.M. Princípios de bioquímica, 5 ed. Porto Alegre: Artmed editora, 2011.`n4. GODFREY, T., WEST, S. (eds), Industrial Enzymology, New York: Chapman-Hall, 1996.`n5. WHITAKER, J.R. (ed.) Pinciples of Enzynmology for the Food Sciences 2nd ed., New York: Marcel Dekker Inc., 1994.`n6. TANAKA, A., TOSA, T., KOBAYASHI, T. (Eds.). Industrial Application of Immobilized Biocatalysts, New York: Marcel Dekker Inc., 1993.`n6.VOET, D., VOET, J., PRATT, C.W. Fundamentos de Bioquímica. Porto Alegre: Editora ARTMED, 2000."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

# --- Row heights / formatting tweaks to match the final layout -------------
$ws.Rows("13:13").RowHeight = 15
$ws.Rows("13:13").EntireRow.AutoFit() | Out-Null

$ws.Rows("14:14").RowHeight = 60
$ws.Rows("15:15").RowHeight = 60
$ws.Rows("16:16").RowHeight = 120
$ws.Rows("17:17").RowHeight = 120
$ws.Rows("18:18").RowHeight = 15
$ws.Rows("22:22").RowHeight = 120

# Narrow column A's explicit range so it no longer overlaps column B's
# width/style definition (column A: 30.71 wide / style 1; column B keeps its
# own 60.71-wide / style 2 definition).
$ws.Columns("A:A").ColumnWidth = 30.7109375
